$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44211
$ws.Range("K2").Value = 'Bing'
$ws.Range("M2").Value = 400
$ws.Range("N2").Value = 6000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 6500
$ws.Range("Q2").Value = '$/caja 7 kilos'
$ws.Range("S2").Value = 929
$ws.Range("T2").Value = 7

# Row 3
$ws.Range("D3").Value = 44211
$ws.Range("M3").Value = 400
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 9500
$ws.Range("Q3").Value = '$/caja 7 kilos'
$ws.Range("S3").Value = 1357
$ws.Range("T3").Value = 7

# Row 4
$ws.Range("D4").Value = 44204
$ws.Range("M4").Value = 400
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 11000
$ws.Range("P4").Value = 10500
$ws.Range("S4").Value = 1050

# Row 5
$ws.Range("D5").Value = 44204
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 15500
$ws.Range("S5").Value = 1550

# Row 6
$ws.Range("D6").Value = 44208
$ws.Range("K6").Value = 'Bing'

# Row 7
$ws.Range("D7").Value = 44208
$ws.Range("K7").Value = 'Rainier'
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 9000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 9500
$ws.Range("Q7").Value = '$/caja 7 kilos'
$ws.Range("R7").Value = 'Provincia de Curicó'
$ws.Range("S7").Value = 1357
$ws.Range("T7").Value = 7

# Row 8
$ws.Range("D8").Value = 44201
$ws.Range("M8").Value = 600
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 11500
$ws.Range("Q8").Value = '$/bandeja 10 kilos'
$ws.Range("S8").Value = 1150
$ws.Range("T8").Value = 10

# Row 9
$ws.Range("D9").Value = 44201
$ws.Range("M9").Value = 600
$ws.Range("N9").Value = 16000
$ws.Range("O9").Value = 17000
$ws.Range("P9").Value = 16500
$ws.Range("Q9").Value = '$/bandeja 10 kilos'
$ws.Range("S9").Value = 1650
$ws.Range("T9").Value = 10

# Row 10
$ws.Range("D10").Value = 44166
$ws.Range("K10").Value = 'Lapins'
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 700
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15500
$ws.Range("Q10").Value = '$/bandeja 10 kilos'
$ws.Range("S10").Value = 1550
$ws.Range("T10").Value = 10

# Row 11
$ws.Range("D11").Value = 44162
$ws.Range("K11").Value = 'Early Burlat'
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 500
$ws.Range("N11").Value = 11000
$ws.Range("O11").Value = 12000
$ws.Range("P11").Value = 11500
$ws.Range("Q11").Value = '$/bandeja 7 kilos'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 1643
$ws.Range("T11").Value = 7

# Row 12
$ws.Range("D12").Value = 44169
$ws.Range("K12").Value = 'Lapins'
$ws.Range("M12").Value = 600
$ws.Range("R12").Value = 'Región de O''Higgins'

# Row 13
$ws.Range("D13").Value = 44189
$ws.Range("M13").Value = 500
$ws.Range("N13").Value = 11000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 11500
$ws.Range("S13").Value = 1150

# Row 14
$ws.Range("D14").Value = 44189
$ws.Range("K14").Value = 'Rainier'
$ws.Range("N14").Value = 16000
$ws.Range("O14").Value = 17000
$ws.Range("P14").Value = 16500
$ws.Range("Q14").Value = '$/bandeja 10 kilos'
$ws.Range("S14").Value = 1650
$ws.Range("T14").Value = 10

# Row 15
$ws.Range("D15").Value = 44159
$ws.Range("K15").Value = 'Early Burlat'
$ws.Range("L15").Value = 'Segunda'
$ws.Range("N15").Value = 11000
$ws.Range("O15").Value = 12000
$ws.Range("P15").Value = 11500
$ws.Range("Q15").Value = '$/bandeja 7 kilos'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 1643
$ws.Range("T15").Value = 7

# Row 16
$ws.Range("D16").Value = 44187
$ws.Range("M16").Value = 600
$ws.Range("N16").Value = 9000
$ws.Range("O16").Value = 10000
$ws.Range("P16").Value = 9500
$ws.Range("Q16").Value = '$/caja 8 kilos'
$ws.Range("S16").Value = 1188
$ws.Range("T16").Value = 8

# Row 17
$ws.Range("D17").Value = 44187
$ws.Range("K17").Value = 'Rainier'
$ws.Range("M17").Value = 600
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 17000
$ws.Range("P17").Value = 16500
$ws.Range("S17").Value = 1650

# Row 18
$ws.Range("D18").Value = 44195
$ws.Range("K18").Value = 'Bing'
$ws.Range("N18").Value = 12000
$ws.Range("O18").Value = 13000
$ws.Range("P18").Value = 12500
$ws.Range("S18").Value = 1250

# Row 19
$ws.Range("D19").Value = 44195
$ws.Range("K19").Value = 'Rainier'
$ws.Range("N19").Value = 16000
$ws.Range("O19").Value = 17000
$ws.Range("P19").Value = 16500
$ws.Range("R19").Value = 'Provincia de Curicó'
$ws.Range("S19").Value = 1650

# Row 20
$ws.Range("D20").Value = 44200

# Row 21
$ws.Range("D21").Value = 44186
$ws.Range("M21").Value = 400
$ws.Range("N21").Value = 15000
$ws.Range("O21").Value = 16000
$ws.Range("P21").Value = 15500
$ws.Range("S21").Value = 1550

# Row 22
$ws.Range("D22").Value = 44217
$ws.Range("K22").Value = 'Bing'
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = 6000
$ws.Range("O22").Value = 7000
$ws.Range("P22").Value = 6500
$ws.Range("Q22").Value = '$/caja 7 kilos'
$ws.Range("R22").Value = 'Provincia de Curicó'
$ws.Range("S22").Value = 929

# Row 23
$ws.Range("D23").Value = 44196
$ws.Range("M23").Value = 200

# Row 24
$ws.Range("D24").Value = 44196
$ws.Range("M24").Value = 200

# Row 25
$ws.Range("D25").Value = 44218
$ws.Range("K25").Value = 'Lapins'
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 300
$ws.Range("N25").Value = 6000
$ws.Range("O25").Value = 7000
$ws.Range("P25").Value = 6500
$ws.Range("Q25").Value = '$/caja 7 kilos'
$ws.Range("R25").Value = 'Provincia de Curicó'
$ws.Range("S25").Value = 929

# Row 26
$ws.Range("D26").Value = 44215
$ws.Range("M26").Value = 300
$ws.Range("N26").Value = 6000
$ws.Range("O26").Value = 7000
$ws.Range("P26").Value = 6500
$ws.Range("Q26").Value = '$/caja 7 kilos'
$ws.Range("S26").Value = 929
$ws.Range("T26").Value = 7

# Row 27
$ws.Range("D27").Value = 44215
$ws.Range("M27").Value = 200
$ws.Range("N27").Value = 9000
$ws.Range("O27").Value = 10000
$ws.Range("P27").Value = 9500
$ws.Range("Q27").Value = '$/caja 7 kilos'
$ws.Range("S27").Value = 1357
$ws.Range("T27").Value = 7

# Row 28
$ws.Range("D28").Value = 44194
$ws.Range("K28").Value = 'Bing'
$ws.Range("M28").Value = 500
$ws.Range("N28").Value = 11000
$ws.Range("O28").Value = 12000
$ws.Range("P28").Value = 11500
$ws.Range("S28").Value = 1150

# Row 29
$ws.Range("D29").Value = 44194
$ws.Range("K29").Value = 'Rainier'
$ws.Range("M29").Value = 600
$ws.Range("N29").Value = 16000
$ws.Range("O29").Value = 17000
$ws.Range("P29").Value = 16500
$ws.Range("Q29").Value = '$/bandeja 10 kilos'
$ws.Range("S29").Value = 1650
$ws.Range("T29").Value = 10

# Row 30
$ws.Range("D30").Value = 44222
$ws.Range("K30").Value = 'Lapins'
$ws.Range("M30").Value = 200
$ws.Range("N30").Value = 6000
$ws.Range("O30").Value = 7000
$ws.Range("P30").Value = 6500
$ws.Range("S30").Value = 929

# Row 31
$ws.Range("D31").Value = 44176
$ws.Range("M31").Value = 500
$ws.Range("N31").Value = 7000
$ws.Range("O31").Value = 8000
$ws.Range("P31").Value = 7500
$ws.Range("Q31").Value = '$/bandeja 5 kilos'
$ws.Range("S31").Value = 1500
$ws.Range("T31").Value = 5

# Row 32
$ws.Range("D32").Value = 44176
$ws.Range("K32").Value = 'Lapins'
$ws.Range("M32").Value = 500
$ws.Range("N32").Value = 13000
$ws.Range("O32").Value = 14000
$ws.Range("P32").Value = 13500
$ws.Range("Q32").Value = '$/bandeja 10 kilos'
$ws.Range("S32").Value = 1350
$ws.Range("T32").Value = 10

# Row 33
$ws.Range("D33").Value = 44161
$ws.Range("K33").Value = 'Early Burlat'
$ws.Range("L33").Value = 'Segunda'
$ws.Range("M33").Value = 400
$ws.Range("N33").Value = 11000
$ws.Range("O33").Value = 12000
$ws.Range("P33").Value = 11500
$ws.Range("Q33").Value = '$/bandeja 7 kilos'
$ws.Range("R33").Value = 'Región de O''Higgins'
$ws.Range("S33").Value = 929
